# Restaurant.xlsx - add "country" column (I) with country data per restaurant group
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("I1").Value = "country"

# British restaurants (London, Brighton, Manchester, Edinburgh, St Albans area)
$ws.Range("I2:I29").Value = "British"
$ws.Range("I30").Value = "French"
$ws.Range("I31").Value = "British"

# French restaurants (Paris, Lyon, etc.)
$ws.Range("I33:I38").Value = "French"
$ws.Range("I40:I45").Value = "French"
$ws.Range("I47:I52").Value = "French"
$ws.Range("I54:I59").Value = "French"
$ws.Range("I61:I66").Value = "French"
$ws.Range("I68:I73").Value = "French"
$ws.Range("I75:I80").Value = "French"

# Spanish restaurants
$ws.Range("I82:I87").Value = "Spain"
$ws.Range("I89:I94").Value = "Spain"
$ws.Range("I96:I101").Value = "Spain"
$ws.Range("I103:I108").Value = "Spain"
$ws.Range("I110:I115").Value = "Spain"

# Italian restaurants
$ws.Range("I117:I122").Value = "Italy"
$ws.Range("I124:I129").Value = "Italy"
$ws.Range("I131:I136").Value = "Italy"
$ws.Range("I138:I143").Value = "Italy"
$ws.Range("I145:I150").Value = "Italy"
$ws.Range("I152:I157").Value = "Italy"
$ws.Range("I159:I164").Value = "Italy"
$ws.Range("I166:I171").Value = "Italy"

# Swiss restaurants
$ws.Range("I173:I178").Value = "Switzerland"
$ws.Range("I180:I185").Value = "Switzerland"
$ws.Range("I187:I192").Value = "Switzerland"
$ws.Range("I194:I199").Value = "Switzerland"
$ws.Range("I201:I206").Value = "Switzerland"
$ws.Range("I208:I213").Value = "Switzerland"
$ws.Range("I215:I220").Value = "Switzerland"
$ws.Range("I222").Value = "Switzerland"

# Match the author's final view state: scrolled so row 181 is at top, with I192 selected
$ws.Range("I192").Select()
try { $excel.ActiveWindow.ScrollRow = 181 } catch {}
